$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 1214
$wsExhibition.Range("F5").Value = 1207
$wsExhibition.Range("G6").Value = 120
$wsExhibition.Range("F7").Value = 4517
$wsExhibition.Range("F8").Value = 2674
$wsExhibition.Range("F10").Value = 2640
$wsExhibition.Range("F15").Value = 697
$wsExhibition.Range("F16").Value = 278
$wsExhibition.Range("F17").Value = 165
$wsExhibition.Range("F18").Value = 353
$wsExhibition.Range("F20").Value = 283
$wsExhibition.Range("F22").Value = 50
$wsExhibition.Range("F26").Value = 593
$wsExhibition.Range("F29").Value = 15
$wsExhibition.Range("F30").Value = 460
$wsExhibition.Range("F31").Value = 1640
$wsExhibition.Range("F32").Value = 1268
$wsExhibition.Range("F33").Value = 239
$wsExhibition.Range("F34").Value = 33
$wsExhibition.Range("F35").Value = 1314
$wsExhibition.Range("F36").Value = 2173
$wsExhibition.Range("F37").Value = 330
$wsExhibition.Range("F39").Value = 570
$wsExhibition.Range("F41").Value = 39
$wsExhibition.Range("F43").Value = 716
$wsExhibition.Range("F44").Value = 1398
$wsExhibition.Range("F45").Value = 153
$wsExhibition.Range("F47").Value = 457
$wsExhibition.Range("F49").Value = 88

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("G2").Value = 88
$wsAllTypes.Range("G3").Value = 88
$wsAllTypes.Range("G4").Value = 120
$wsAllTypes.Range("F5").Value = 4517
$wsAllTypes.Range("F6").Value = 2674
$wsAllTypes.Range("F7").Value = 2640
$wsAllTypes.Range("F11").Value = 697
$wsAllTypes.Range("F12").Value = 278
$wsAllTypes.Range("F13").Value = 165
$wsAllTypes.Range("F14").Value = 353
$wsAllTypes.Range("F16").Value = 283
$wsAllTypes.Range("F18").Value = 50
$wsAllTypes.Range("F21").Value = 593
$wsAllTypes.Range("F27").Value = 460
$wsAllTypes.Range("F28").Value = 1640
$wsAllTypes.Range("F29").Value = 1268
$wsAllTypes.Range("F30").Value = 239
$wsAllTypes.Range("F31").Value = 33
$wsAllTypes.Range("F34").Value = 2173
$wsAllTypes.Range("F35").Value = 330
$wsAllTypes.Range("F39").Value = 570
$wsAllTypes.Range("F41").Value = 39
$wsAllTypes.Range("F43").Value = 716
$wsAllTypes.Range("F44").Value = 1398
$wsAllTypes.Range("F46").Value = 153
$wsAllTypes.Range("F47").Value = 457
$wsAllTypes.Range("F48").Value = 88
